$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D is treated as text for the cells we touch, to preserve exact string formatting

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Simple Price-only (column D) updates ---
Set-TextValue $ws.Range("D2") "272.10"
Set-TextValue $ws.Range("D3") "23.00"
Set-TextValue $ws.Range("D4") "6.379"
Set-TextValue $ws.Range("D5") "0.06255"
Set-TextValue $ws.Range("D7") "6.763"
Set-TextValue $ws.Range("D8") "1.399"
Set-TextValue $ws.Range("D9") "0.8384"
Set-TextValue $ws.Range("D10") "0.01372"
Set-TextValue $ws.Range("D11") "0.1621"
Set-TextValue $ws.Range("D12") "0.08349"
Set-TextValue $ws.Range("D13") "0.03414"
Set-TextValue $ws.Range("D14") "0.03181"
Set-TextValue $ws.Range("D40") "0.04683"
Set-TextValue $ws.Range("D41") "0.006915"
Set-TextValue $ws.Range("D42") "0.1168"
Set-TextValue $ws.Range("D43") "0.003481"
Set-TextValue $ws.Range("D44") "0.01250"
Set-TextValue $ws.Range("D45") "0.00006305"
Set-TextValue $ws.Range("D46") "0.00000000754"
Set-TextValue $ws.Range("D47") "0.7044"
Set-TextValue $ws.Range("D48") "0.1282"
Set-TextValue $ws.Range("D49") "0.00002113"
Set-TextValue $ws.Range("D50") "0.01247"

# --- Rows 15-26: coin list shifted down by one (ProBitToken moved to top), with updated prices ---
$ws.Range("B15").Value = "ProBitToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws.Range("D15") "0.1253"
$ws.Range("E15").Value = "14ProBitTokenPROB"

$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D16") "0.09307"
$ws.Range("E16").Value = "15BitMartTokenBMX"

$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D17") "3.925"
$ws.Range("E17").Value = "16MCDexMCB"

$ws.Range("B18").Value = "BitForexToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D18") "0.001708"
$ws.Range("E18").Value = "17BitForexTokenBF"

$ws.Range("B19").Value = "CoinExToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D19") "0.04909"
$ws.Range("E19").Value = "18CoinExTokenCET"

$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D20") "0.006274"
$ws.Range("E20").Value = "19TigerCashTCH"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D21") "0.005432"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D22") "0.001097"
$ws.Range("E22").Value = "21BitKanKAN"

$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D23") "0.0001509"
$ws.Range("E23").Value = "22NitroExNTX"

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D24") "3.732"
$ws.Range("E24").Value = "23LEOLEO"

$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D25") "2.330"
$ws.Range("E25").Value = "24BTSETokenBTSE"

$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Range("D26") "0.3338"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
